# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500adfae01c9a5dd7ef65e90accc96781b5c
#
# Re-brand the IBM/Alvearie published artifact as a LinuxForHealth one:
#  - Metadata sheet: update URL, Version, Date and Publisher
#  - Elements sheet: the "Constraint(s)" text for the Extension row was a
#    stray duplicate of the constraint that belongs to Extension.extension;
#    clear it from the Extension row (row 2) so it only appears once.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/usual-and-customary-amount"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
